# Applies the edits described by the commit "letshope this is almost done"
# to the RBD standardized questionnaire workbook.

$wb = $excel.ActiveWorkbook

$survey  = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# -----------------------------------------------------------------
# 1) "survey" sheet: stamp column G ("yes") on a bunch of rows that
#    previously had no relevance/choice_filter value in that column.
# -----------------------------------------------------------------
$surveyGRows = @(40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,58,59,102,103,104,105,106,107,108,109,114)
foreach ($r in $surveyGRows) {
    $survey.Cells.Item($r, 7).Value = "yes"
}

# -----------------------------------------------------------------
# 2) "choices" sheet: add two new list choices (new shared strings).
#    Order matters for shared-string index assignment: A127 first so
#    "bonn+A5:C127espratiques" lands before "Yesno".
# -----------------------------------------------------------------
$choices.Range("A127").Value = "bonn+A5:C127espratiques"
$choices.Range("A24").Value = "Yesno"
$choices.Range("A25").Value = "Yesno"

# -----------------------------------------------------------------
# 3) "choices" sheet: numeric tweaks on a handful of rows.
# -----------------------------------------------------------------
$choices.Range("B66").Value = 888
$choices.Range("B67").Value = 8888
$choices.Range("B72").Value = 888
$choices.Range("B73").Value = 8888
$choices.Range("B81").Value = 888
$choices.Range("B82").Value = 8888
$choices.Range("B102").Value = 888
$choices.Range("B103").Value = 8888

# -----------------------------------------------------------------
# 4) Selection / active-sheet bookkeeping to mirror the saved view
#    state recorded in the diff (choices becomes the active/visible
#    tab, survey's selection moves to G114, choices' selection moves
#    to A25).
# -----------------------------------------------------------------
$survey.Activate()
$survey.Range("G114").Select()

$choices.Activate()
$choices.Range("A25").Select()

Write-Output "edits applied"
